# Auto-generated edit script: updates market-data derived columns (H:N)
# on multiple worksheets to reflect refreshed pricing figures from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 722.75
$ws.Range("I6").Value2 = 510.5
$ws.Range("J6").Value2 = 935
$ws.Range("K6").Value2 = 1531.5
$ws.Range("L6").Value2 = 2805
$ws.Range("M6").Value2 = -1419.5
$ws.Range("N6").Value2 = -3029
$ws.Range("H12").Value2 = 148.4
$ws.Range("I12").Value2 = 100
$ws.Range("K12").Value2 = 100
$ws.Range("M12").Value2 = 70
$ws.Range("H17").Value2 = 1692.8
$ws.Range("J17").Value2 = 1878.375
$ws.Range("L17").Value2 = 5635.125
$ws.Range("N17").Value2 = -5971.125
$ws.Range("H21").Value2 = 18.6
$ws.Range("J21").Value2 = 18.6
$ws.Range("L21").Value2 = 18.6
$ws.Range("N21").Value2 = -954.6
$ws.Range("H23").Value2 = 18.6
$ws.Range("J23").Value2 = 18.6
$ws.Range("L23").Value2 = 18.6
$ws.Range("N23").Value2 = -486.6
$ws.Range("H29").Value2 = 1474.5
$ws.Range("J29").Value2 = 0
$ws.Range("L29").Value2 = 0
$ws.Range("N29").ClearContents() | Out-Null
$ws.Range("H115").Value2 = 799.1429000000001
$ws.Range("I115").Value2 = 706.7692
$ws.Range("K115").Value2 = 2120.3076
$ws.Range("M115").Value2 = -553.3076000000001
$ws.Range("H125").Value2 = 7410443.5
$ws.Range("I125").Value2 = 1250
$ws.Range("K125").Value2 = 11250
$ws.Range("M125").Value2 = -8790
$ws.Range("H138").Value2 = 6595.8735
$ws.Range("J138").Value2 = 6942.258
$ws.Range("L138").Value2 = 20826.774
$ws.Range("N138").Value2 = -31106.774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 348654.25
$ws.Range("I74").Value2 = 458180.6
$ws.Range("K74").Value2 = 458180.6
$ws.Range("M74").Value2 = -457306.6
$ws.Range("H77").Value2 = 348654.25
$ws.Range("I77").Value2 = 458180.6
$ws.Range("K77").Value2 = 2290903
$ws.Range("M77").Value2 = -2286535
$ws.Range("H102").Value2 = 3105.5454
$ws.Range("I102").Value2 = 3025
$ws.Range("J102").Value2 = 3911
$ws.Range("K102").Value2 = 3025
$ws.Range("L102").Value2 = 3911
$ws.Range("M102").Value2 = -1403
$ws.Range("N102").Value2 = -7155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 2095
$ws.Range("I20").Value2 = 2174.2
$ws.Range("K20").Value2 = 2174.2
$ws.Range("M20").Value2 = -1927.2
$ws.Range("H81").Value2 = 16106
$ws.Range("J81").Value2 = 16106
$ws.Range("L81").Value2 = 16106
$ws.Range("N81").Value2 = -18228
$ws.Range("H84").Value2 = 16106
$ws.Range("J84").Value2 = 16106
$ws.Range("L84").Value2 = 48318
$ws.Range("N84").Value2 = -58926
$ws.Range("H134").Value2 = 279379.44
$ws.Range("I134").Value2 = 307444.8
$ws.Range("J134").Value2 = 147071.28
$ws.Range("K134").Value2 = 922334.3999999999
$ws.Range("L134").Value2 = 441213.84
$ws.Range("M134").Value2 = -919799.3999999999
$ws.Range("N134").Value2 = -446283.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 90582.25
$ws.Range("I99").Value2 = 8599.799999999999
$ws.Range("K99").Value2 = 8599.799999999999
$ws.Range("M99").Value2 = -7101.799999999999
$ws.Range("H105").Value2 = 0
$ws.Range("I105").Value2 = 0
$ws.Range("J105").Value2 = 0
$ws.Range("K105").Value2 = 0
$ws.Range("L105").Value2 = 0
$ws.Range("M105").ClearContents() | Out-Null
$ws.Range("N105").ClearContents() | Out-Null
$ws.Range("H126").Value2 = 90582.25
$ws.Range("I126").Value2 = 8599.799999999999
$ws.Range("K126").Value2 = 25799.4
$ws.Range("M126").Value2 = -23329.4
$ws.Range("H132").Value2 = 2903
$ws.Range("I132").Value2 = 2462.3572
$ws.Range("J132").Value2 = 5987.5
$ws.Range("K132").Value2 = 7387.071599999999
$ws.Range("L132").Value2 = 17962.5
$ws.Range("M132").Value2 = -4857.071599999999
$ws.Range("N132").Value2 = -23022.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 2346.1304
$ws.Range("J34").Value2 = 2638.1
$ws.Range("L34").Value2 = 7914.299999999999
$ws.Range("N34").Value2 = -8082.299999999999
$ws.Range("H60").Value2 = 745.4
$ws.Range("H121").Value2 = 1748.5
$ws.Range("I121").Value2 = 0
$ws.Range("J121").Value2 = 1748.5
$ws.Range("K121").Value2 = 0
$ws.Range("L121").Value2 = 5245.5
$ws.Range("M121").ClearContents() | Out-Null
$ws.Range("N121").Value2 = -7865.5
$ws.Range("H131").Value2 = 29924.576
$ws.Range("I131").Value2 = 73172.5
$ws.Range("J131").Value2 = 19662.355
$ws.Range("K131").Value2 = 219517.5
$ws.Range("L131").Value2 = 58987.065
$ws.Range("M131").Value2 = -214477.5
$ws.Range("N131").Value2 = -69067.065
$ws.Range("H137").Value2 = 4343.5
$ws.Range("I137").Value2 = 2678.7307
$ws.Range("J137").Value2 = 9754
$ws.Range("K137").Value2 = 8036.1921
$ws.Range("L137").Value2 = 29262
$ws.Range("M137").Value2 = -2936.1921
$ws.Range("N137").Value2 = -39462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 1346261.8
$ws.Range("I80").Value2 = 1678495.5
$ws.Range("J80").Value2 = 1124772.5
$ws.Range("K80").Value2 = 1678495.5
$ws.Range("L80").Value2 = 1124772.5
$ws.Range("M80").Value2 = -1677497.5
$ws.Range("N80").Value2 = -1126768.5
$ws.Range("H83").Value2 = 1346261.8
$ws.Range("I83").Value2 = 1678495.5
$ws.Range("J83").Value2 = 1124772.5
$ws.Range("K83").Value2 = 8392477.5
$ws.Range("L83").Value2 = 5623862.5
$ws.Range("M83").Value2 = -8387485.5
$ws.Range("N83").Value2 = -5633846.5
$ws.Range("H132").Value2 = 1240124.1
$ws.Range("I132").Value2 = 1270139.6
$ws.Range("K132").Value2 = 3810418.8
$ws.Range("M132").Value2 = -3807888.8
$ws.Range("H134").Value2 = 49000
$ws.Range("J134").Value2 = 49000
$ws.Range("L134").Value2 = 147000
$ws.Range("N134").Value2 = -152070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 6789.5
$ws.Range("I7").Value2 = 5125.25
$ws.Range("J7").Value2 = 7899
$ws.Range("K7").Value2 = 5125.25
$ws.Range("L7").Value2 = 7899
$ws.Range("M7").Value2 = -5013.25
$ws.Range("N7").Value2 = -8123
$ws.Range("H20").Value2 = 1260600
$ws.Range("J20").Value2 = 1568750
$ws.Range("L20").Value2 = 1568750
$ws.Range("N20").Value2 = -1569202
$ws.Range("H22").Value2 = 800
$ws.Range("I22").Value2 = 800
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 800
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -505
$ws.Range("N22").ClearContents() | Out-Null
$ws.Range("H27").Value2 = 800
$ws.Range("I27").Value2 = 800
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 800
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = -693
$ws.Range("N27").ClearContents() | Out-Null
$ws.Range("H43").Value2 = 16453000
$ws.Range("I43").Value2 = 382500
$ws.Range("K43").Value2 = 382500
$ws.Range("M43").Value2 = -382307
$ws.Range("H46").Value2 = 2793.9443
$ws.Range("I46").Value2 = 2823
$ws.Range("J46").Value2 = 2718.4
$ws.Range("K46").Value2 = 2823
$ws.Range("L46").Value2 = 2718.4
$ws.Range("M46").Value2 = -2635
$ws.Range("N46").Value2 = -3094.4
$ws.Range("H68").Value2 = 4496.3335
$ws.Range("I68").Value2 = 3919
$ws.Range("J68").Value2 = 4958.2
$ws.Range("K68").Value2 = 3919
$ws.Range("L68").Value2 = 4958.2
$ws.Range("M68").Value2 = -3170
$ws.Range("N68").Value2 = -6456.2
$ws.Range("H71").Value2 = 4496.3335
$ws.Range("I71").Value2 = 3919
$ws.Range("J71").Value2 = 4958.2
$ws.Range("K71").Value2 = 19595
$ws.Range("L71").Value2 = 24791
$ws.Range("M71").Value2 = -15851
$ws.Range("N71").Value2 = -32279
$ws.Range("H82").Value2 = 1660.0526
$ws.Range("I82").Value2 = 1691
$ws.Range("J82").Value2 = 1637.5454
$ws.Range("K82").Value2 = 1691
$ws.Range("L82").Value2 = 1637.5454
$ws.Range("M82").Value2 = -1330
$ws.Range("N82").Value2 = -2359.5454
$ws.Range("H85").Value2 = 1660.0526
$ws.Range("I85").Value2 = 1691
$ws.Range("J85").Value2 = 1637.5454
$ws.Range("K85").Value2 = 1691
$ws.Range("L85").Value2 = 1637.5454
$ws.Range("M85").Value2 = -443
$ws.Range("N85").Value2 = -4133.5454
$ws.Range("H112").Value2 = 0
$ws.Range("J112").Value2 = 0
$ws.Range("L112").Value2 = 0
$ws.Range("N112").ClearContents() | Out-Null
$ws.Range("H126").Value2 = 6789.5
$ws.Range("I126").Value2 = 5125.25
$ws.Range("J126").Value2 = 7899
$ws.Range("K126").Value2 = 15375.75
$ws.Range("L126").Value2 = 23697
$ws.Range("M126").Value2 = -12905.75
$ws.Range("N126").Value2 = -28637
$ws.Range("H132").Value2 = 5069.7427
$ws.Range("I132").Value2 = 4447.577
$ws.Range("K132").Value2 = 13342.731
$ws.Range("M132").Value2 = -10812.731
$ws.Range("H136").Value2 = 563322.7
$ws.Range("I136").Value2 = 774754.4399999999
$ws.Range("J136").Value2 = 13600
$ws.Range("K136").Value2 = 2324263.32
$ws.Range("L136").Value2 = 40800
$ws.Range("M136").Value2 = -2321713.32
$ws.Range("N136").Value2 = -45900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 8212.454
$ws.Range("I62").Value2 = 10334.333
$ws.Range("J62").Value2 = 7416.75
$ws.Range("K62").Value2 = 10334.333
$ws.Range("L62").Value2 = 7416.75
$ws.Range("M62").Value2 = -9710.333000000001
$ws.Range("N62").Value2 = -8664.75
$ws.Range("H65").Value2 = 8212.454
$ws.Range("I65").Value2 = 10334.333
$ws.Range("J65").Value2 = 7416.75
$ws.Range("K65").Value2 = 51671.665
$ws.Range("L65").Value2 = 37083.75
$ws.Range("M65").Value2 = -48551.665
$ws.Range("N65").Value2 = -43323.75
$ws.Range("H81").Value2 = 11222.111
$ws.Range("I81").Value2 = 4374.875
$ws.Range("K81").Value2 = 8749.75
$ws.Range("M81").Value2 = -7688.75
$ws.Range("H84").Value2 = 11222.111
$ws.Range("I84").Value2 = 4374.875
$ws.Range("K84").Value2 = 43748.75
$ws.Range("M84").Value2 = -38444.75
$ws.Range("H132").Value2 = 69403.8
$ws.Range("I132").Value2 = 2932.6428
$ws.Range("J132").Value2 = 1000000
$ws.Range("K132").Value2 = 8797.928400000001
$ws.Range("L132").Value2 = 3000000
$ws.Range("M132").Value2 = -6267.928400000001
$ws.Range("N132").Value2 = -3005060
